$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Generated:" timestamp footer cell
$ws.Range("A30").Value = "Generated: 2023-09-04 10:47:14 AM"

# Row 15 (Naveen Bromiyo A R) is now Completed:
#  - Student Name (B15) and Completion Status (E15) get the "Completed" look
#    (bold white font on green fill), matching the styling already used by
#    every other Completed row (e.g. row 2).
#  - Completion Status (E15) changes from "Pending" to "Completed".
#  - Pending Task (D15) is cleared since the task is finished.
$ws.Range("E15").Value = "Completed"

$ws.Range("B2").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("D15").ClearContents()
